$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.304.54"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.02%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.087.12"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.40%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.65%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "342.59"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.84%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.61%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5228"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +1.54%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4412"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.21%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "54.31"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +3.35%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09315"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.27%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.167"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.47%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.77"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.44%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.577"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +3.56%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.888"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.94%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.089.66"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.73%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "101.05"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.44%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001159"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.64%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.002"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.71%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "21.08"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.11%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.06654"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.14%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.325"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.95%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.43%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.345.30"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.49"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.97%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.300"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.59%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "21.76"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.75%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "162.43"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.28%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.511"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.59%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "132.88"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.16%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.133"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.06%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.659"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.39%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.62%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.222"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.79%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.672"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +8.32%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.854"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.12"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.17%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02623"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.84%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06830"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.63%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +3.82%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6956"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +1.48%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.49"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.29%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.2200"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.06%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6806"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +2.66%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.37"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.29%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.322"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.14%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.48%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.372"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +18.26%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.632"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.14%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00000000344"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.32%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.206"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +7.52%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.212"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.68%  "
